$d = $word.ActiveDocument

$replacements = @(
    @("397×4=", "870×6="),
    @("275×2=", "237×6="),
    @("361×4=", "356×8="),
    @("922×9=", "995×2="),
    @("316×7=", "534×5="),
    @("735×6=", "329×7="),
    @("583×8=", "196×5="),
    @("516×2=", "949×9="),
    @("915×3=", "921×4="),
    @("271×5=", "459×2="),
    @("403×2=", "242×9="),
    @("905×3=", "749×2="),
    @("197×3=", "692×7="),
    @("735×4=", "707×8="),
    @("516×6=", "315×8="),
    @("780×9=", "515×6="),
    @("182×2=", "398×3="),
    @("356×2=", "325×7="),
    @("424×4=", "465×2="),
    @("581×5=", "353×8="),
    @("840×5=", "647×2="),
    @("235×6=", "963×2="),
    @("593×5=", "422×4="),
    @("138×8=", "645×9="),
    @("615×9=", "362×6="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
